$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.861.41'
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').Value = '2.292.18'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +0.65%  '
$ws.Range('D5').Value = '314.64'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('D6').Value = '103.36'
$ws.Range('E6').Value = '  +0.76%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('D9').Value = '0.602'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('D10').Value = '39.25'
$ws.Range('E10').Value = '  -0.94%  '
$ws.Range('D11').Value = '0.0906'
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').Value = '8.33'
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('E13').Value = '  +2.49%  '
$ws.Range('E14').Value = '  +2.77%  '
$ws.Range('D15').Value = '15.20'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').Value = '2.639.80'
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('D17').Value = '2.288.73'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').Value = '42.781.99'
$ws.Range('E18').Value = '  +0.99%  '
$ws.Range('D19').Value = '7.42'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').Value = '13.69'
$ws.Range('E20').Value = '  +17.15%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0000105'
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('D22').Value = '73.71'
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('E23').Value = '  +0.76%  '
$ws.Range('D24').Value = '265.17'
$ws.Range('E24').Value = '  -3.50%  '
$ws.Range('E25').Value = '  -2.63%  '
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('D27').Value = '10.78'
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('D28').Value = '2.34'
$ws.Range('E28').Value = '  -0.86%  '
$ws.Range('D29').Value = '7.00'
$ws.Range('E29').Value = '  +17.72%  '
$ws.Range('D30').Value = '22.51'
$ws.Range('E30').Value = '  -0.86%  '
$ws.Range('D31').Value = '36.65'
$ws.Range('E31').Value = '  -2.00%  '
$ws.Range('D32').Value = '167.18'
$ws.Range('E32').Value = '  +0.89%  '
$ws.Range('D33').Value = '0.0868'
$ws.Range('E33').Value = '  -0.66%  '
$ws.Range('E34').Value = '  -2.25%  '
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('E36').Value = '  -5.00%  '
$ws.Range('D37').Value = '4.52'
$ws.Range('E37').Value = '  -1.12%  '
$ws.Range('E38').Value = '  -3.65%  '
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('E40').Value = '  -3.48%  '
$ws.Range('E41').Value = '  +5.71%  '
$ws.Range('D42').Value = '70.06'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('E43').Value = '  +2.00%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('B45').Value = 'BitcoinSV'
$ws.Range('C45').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D45').Value = '93.84'
$ws.Range('E45').Value = '  -0.96%  '
$ws.Range('D46').Value = '1.735.80'
$ws.Range('E46').Value = '  +9.25%  '
$ws.Range('D47').Value = '11.99'
$ws.Range('D48').Value = '79.87'
$ws.Range('E48').Value = '  -0.88%  '
$ws.Range('D49').Value = '111.92'
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('E50').Value = '  -0.53%  '
$ws.Range('E51').Value = '  -3.06%  '
